$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A208").Value = "25-10-2021"
$ws.Range("B208").Value = 5.18
$ws.Range("C208").Value = 5.3
$ws.Range("D208").Value = 5.43
$ws.Range("E208").Value = 5.54
$ws.Range("F208").Value = 5.82
$ws.Range("G208").Value = -0.68
$ws.Range("H208").Value = 0.38
$ws.Range("I208").Value = 1.04
$ws.Range("J208").Value = 1.41
$ws.Range("K208").Value = 1.81
$ws.Range("L208").Value = 2.39
$ws.Range("M208").Value = 2.44

$ws.Range("A209").Value = "26-10-2021"
$ws.Range("B209").Value = 5.18
$ws.Range("C209").Value = 5.29
$ws.Range("D209").Value = 5.41
$ws.Range("E209").Value = 5.52
$ws.Range("F209").Value = 5.78
$ws.Range("G209").Value = -0.83
$ws.Range("H209").Value = 0.31
$ws.Range("I209").Value = 0.98
$ws.Range("J209").Value = 1.34
$ws.Range("K209").Value = 1.73
$ws.Range("L209").Value = 2.28
$ws.Range("M209").Value = 2.36

$ws.Range("A210").Value = "27-10-2021"
$ws.Range("B210").Value = 5.15
$ws.Range("C210").Value = 5.24
$ws.Range("D210").Value = 5.36
$ws.Range("E210").Value = 5.46
$ws.Range("F210").Value = 5.73
$ws.Range("G210").Value = -0.88
$ws.Range("H210").Value = 0.27
$ws.Range("I210").Value = 0.96
$ws.Range("J210").Value = 1.32
$ws.Range("K210").Value = 1.67
$ws.Range("L210").Value = 2.2
$ws.Range("M210").Value = 2.26

$ws.Range("A211").Value = "28-10-2021"
$ws.Range("B211").Value = 5.2
$ws.Range("C211").Value = 5.25
$ws.Range("D211").Value = 5.34
$ws.Range("E211").Value = 5.43
$ws.Range("F211").Value = 5.63
$ws.Range("G211").Value = -0.77
$ws.Range("H211").Value = 0.25
$ws.Range("I211").Value = 0.97
$ws.Range("J211").Value = 1.29
$ws.Range("K211").Value = 1.56
$ws.Range("L211").Value = 2.07
$ws.Range("M211").Value = 2.13

$ws.Range("A212").Value = "29-10-2021"
$ws.Range("B212").Value = 5.26
$ws.Range("C212").Value = 5.3
$ws.Range("D212").Value = 5.38
$ws.Range("E212").Value = 5.47
$ws.Range("F212").Value = 5.68
$ws.Range("G212").Value = -1.01
$ws.Range("H212").Value = 0.21
$ws.Range("I212").Value = 0.93
$ws.Range("J212").Value = 1.27
$ws.Range("K212").Value = 1.51
$ws.Range("L212").Value = 2.04
$ws.Range("M212").Value = 2.1
